$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 209 ("「ラヤンはねむれない」" entry); all rows below shift up by one.
$ws.Rows.Item(209).Delete()
